$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Update the cover subtitle text: version (530) -> (561)
$ws.Range("B2").Value = "BSI Super Apps - App Version 1.0.2 (561) Rebrand OCP QA"

# Restore the normal single-cell selection on the active sheet (was a full-column selection)
$ws.Activate()
$ws.Range("B3").Select()
